$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F ("想去人数") values for these rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 436
$ws1.Range("F9").Value = 6995
$ws1.Range("F13").Value = 7893
$ws1.Range("F16").Value = 5467
$ws1.Range("F17").Value = 46
$ws1.Range("F18").Value = 2351
$ws1.Range("F19").Value = 1001
$ws1.Range("F20").Value = 4543
$ws1.Range("F21").Value = 281
$ws1.Range("F22").Value = 375
$ws1.Range("F25").Value = 337
$ws1.Range("F28").Value = 2173
$ws1.Range("F30").Value = 250
$ws1.Range("F32").Value = 83
$ws1.Range("F33").Value = 557
$ws1.Range("F35").Value = 25
$ws1.Range("F36").Value = 1440
$ws1.Range("F39").Value = 2198
$ws1.Range("F40").Value = 2194

# Sheet "演出" (sheet2): update column F ("想去人数") value for this row
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 86

# Sheet "全部类型" (sheet4): update column F ("想去人数") values for these rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 86
$ws4.Range("F11").Value = 436
$ws4.Range("F12").Value = 6995
$ws4.Range("F16").Value = 7893
$ws4.Range("F19").Value = 5467
$ws4.Range("F20").Value = 46
$ws4.Range("F21").Value = 2351
$ws4.Range("F22").Value = 1001
$ws4.Range("F23").Value = 4543
$ws4.Range("F24").Value = 281
$ws4.Range("F25").Value = 375
$ws4.Range("F30").Value = 337
$ws4.Range("F33").Value = 2173
$ws4.Range("F35").Value = 250
$ws4.Range("F37").Value = 83
$ws4.Range("F38").Value = 557
$ws4.Range("F40").Value = 25
$ws4.Range("F42").Value = 1440
$ws4.Range("F45").Value = 2198
$ws4.Range("F47").Value = 2194
